# convection.xlsx - "Add files via upload"
#
# The original commit adds two new quiz sheets ("9_" and "10_") after the
# existing "8_" sheet, renumbers/extends the shared-string table, tweaks the
# existing "8_" sheet's header question + row height + selection, and moves
# the active/selected tab to the new last sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the existing "8_" sheet (steel-ball cooling-curve question):
#    - reword the header question in A1
#    - row 1 grows from 120 to 150 (longer wrapped text)
#    - selection moves off this sheet once new sheets are added
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("8_")

$ws8.Range("A1").Value = "If the starting temperature of the steel ball is 100 degrees C, and the liquid is held constant at 0 degrees C, what would we expect the plot of the temperature over time to look like?  Of the plots in the image above seems the most likely portrait of the temperature behavior of the sphere?"
$ws8.Rows(1).RowHeight = 150

# ---------------------------------------------------------------------
# 2. Add sheet "9_" right after "8_" - derive it from "8_" so it inherits
#    the same quiz-style formatting (header fill, column widths, etc.)
# ---------------------------------------------------------------------
$ws8.Copy($null, $ws8)
$ws9 = $wb.ActiveSheet
$ws9.Name = "9_"

$ws9.Range("A1").Value = "Look at the equation for temperature that we just derived.  What is the temperature of the object at time 't = 0', according to the equation?"
$ws9.Rows(1).RowHeight = 90

$ws9.Range("A2").Value = 0
$ws9.Range("B2").Value = "N"
$ws9.Range("C2").ClearContents()

$ws9.Range("A3").Value = 1
$ws9.Range("B3").Value = "N"
$ws9.Range("C3").ClearContents()

$ws9.Range("A4").Value = "T_infinity"
$ws9.Range("B4").Value = "N"
$ws9.Range("C4").ClearContents()
$ws9.Rows(4).RowHeight = 15

$ws9.Range("A5").Value = "T_i "
$ws9.Range("B5").Value = "Y"
$ws9.Range("C5").Value = "Just as you'd expect.  E^0 is 1, and this means the 2 'T_infinity' terms cancel out, leaving only 'T_i'"
$ws9.Rows(5).RowHeight = 45

$ws9.Range("A6").Value = -1
$ws9.Range("B6").Value = "N"
$ws9.Range("C6").ClearContents()

$ws9.Range("A1:C6").Select()

# ---------------------------------------------------------------------
# 3. Add sheet "10_" right after "9_"
# ---------------------------------------------------------------------
$ws9.Copy($null, $ws9)
$ws10 = $wb.ActiveSheet
$ws10.Name = "10_"

$ws10.Range("A1").Value = "Look at the equation for temperature that we just derived.  What is the temperature of the object at time 't = infinity', according to the equation?"
$ws10.Rows(1).RowHeight = 90

$ws10.Range("A2").Value = 0
$ws10.Range("B2").Value = "N"
$ws10.Range("C2").ClearContents()

$ws10.Range("A3").Value = 1
$ws10.Range("B3").Value = "N"
$ws10.Range("C3").ClearContents()

$ws10.Range("A4").Value = "T_infinity"
$ws10.Range("B4").Value = "Y"
$ws10.Range("C4").Value = "e^(- infinity) = 0, so the only term left is 'T_infinity'.  This makes sense physically: as the steel sphere cools, it will slowly approach the value of the liquid, 'T_infinity'"
$ws10.Rows(4).RowHeight = 75

$ws10.Range("A5").Value = "T_i "
$ws10.Range("B5").Value = "N"
$ws10.Range("C5").ClearContents()
$ws10.Rows(5).RowHeight = 15

$ws10.Range("A6").Value = -1
$ws10.Range("B6").Value = "N"
$ws10.Range("C6").ClearContents()

$ws10.Range("D16").Select()
$ws10.Select()
